# TC27_Canine_Filter_Breed-Giant.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The "CasesTab" Cypher query stored in B2 referenced an unused/erroring
# `co` (cohort) variable in its final RETURN line
# (`coalesce(co.cohort_description, '') AS `Cohort``). That trailing
# line (and the blank line after it) is removed here; the `OPTIONAL
# MATCH (co:cohort)...` / `WITH ... co` plumbing is left alone since the
# query still matches it, it's just no longer returned.
#
# B3 (SamplesTab) / B4 (FilesTab) queries are unaffected content-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Giant Schnauzer'] `nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value2 = $casesQuery

# Row heights re-wrap shorter now that the Cases query lost two lines;
# Samples/Files rows keep the same line count but shift by the new
# (slightly smaller) default line height.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# Window/selection state left by whoever last saved the workbook.
$win = $excel.ActiveWindow
$win.Zoom = 85
$ws.Range("B2").Select()

# Best-effort: reflect the Excel application window geometry recorded in
# the saved file (not all hosts surface this back into the OOXML, but
# this is the correct COM call for it).
$win2 = $excel.Windows.Item(1)
$win2.Left = 28680
$win2.Top = -105
$win2.Width = 29040
$win2.Height = 15840
